$wb = $excel.ActiveWorkbook

# --- Sheet "SoCDTtiNTY-psgr" (passenger) ---
$wsPsgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")

# Row 2 (LDVs): 0.0755 -> 0.076 across B2:H2
$wsPsgr.Range("B2:H2").Value = 0.076

# Row 7 (motorbikes): 0.0587 -> 0.0582 across B7:H7
$wsPsgr.Range("B7:H7").Value = 0.0582

# --- Sheet "SoCDTtiNTY-frgt" (freight) ---
$wsFrgt = $wb.Worksheets.Item("SoCDTtiNTY-frgt")

# Row 2 (LDVs): 0.07 -> 0.078 (B2) and 0.0755 (C2:H2)
$wsFrgt.Range("B2").Value = 0.078
$wsFrgt.Range("C2:H2").Value = 0.0755

# --- View / selection state ---
# Select the freight sheet's range first, then finally activate the
# passenger sheet so it ends up as the selected/active tab, matching
# the target workbook view state.
$wsFrgt.Activate()
$wsFrgt.Range("B2:H2").Select()

$wsPsgr.Activate()
$wsPsgr.Range("B2:H7").Select()
